$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing "医务人员满意度▲" / "医务人员满意度"
# (row 96), shifting all subsequent rows up by one.
$ws.Rows.Item(96).Delete()

# Match the new scroll position / active cell left behind by the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 88
$win.ScrollColumn = 1
$ws.Range("A96").Select()
